# Loan RBI, Variable Instalments
# The "Repayment schedule" worksheet gains a new (blank) column inserted
# before the former "Late" column (column N), shifting the former
# N/O/P columns (Late / heading-spacer / Outstanding) one position to
# the right (O/P/Q).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)   # "Repayment schedule"

# Insert a new blank column at position N (14), pushing existing
# N:P columns to O:Q.
$ws.Columns.Item(14).Insert() | Out-Null

# The newly inserted column keeps the same width as its neighbour
# (column M / 13) rather than the autosized "bestFit" width of the
# column that used to live there.
$ws.Columns.Item(14).ColumnWidth = 10.14

# Restore the last selected cell as stored in the saved workbook.
$ws.Range("I16").Select() | Out-Null
